{"js": "// Center the second table (the stock-data table) and round its\n// \"Cot.M\u00e1xima R$\" column (index 2) values to two decimal places.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length < 2) {\n  throw new Error(\"Expected at least 2 tables in the document\");\n}\n\n// The data table (with the numeric rows) is the second table on the page.\nconst dataTable = tables.items[1];\n\n// 1) Center the table (<w:jc w:val=\"center\"/> on <w:tblPr>).\ndataTable.alignment = \"Centered\";\n\n// 2) Round the \"Cot.M\u00e1xima R$\" column (3rd column, index 2) values.\ndataTable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = dataTable.rows.items;\nfor (const row of rows) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nconst targetCells = rows.map((row) => row.cells.items[2]);\nfor (const cell of targetCells) {\n  cell.load(\"value\");\n}\nawait context.sync();\n\nfor (const cell of targetCells) {\n  const current = cell.value;\n  const num = parseFloat(current);\n  if (!isNaN(num)) {\n    const rounded = Math.round(num * 100) / 100;\n    cell.value = String(rounded);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Center the second table (the stock-data table) and round its\n# \"Cot.M\u00e1xima R$\" column (3rd column) values to two decimal places.\n\n$d = $word.ActiveDocument\n\nif ($d.Tables.Count -lt 2) {\n    throw \"Expected at least 2 tables in the document\"\n}\n\n# The data table (with the numeric rows) is the second table on the page.\n$t = $d.Tables.Item(2)\n\n# 1) Center the table (<w:jc w:val=\"center\"/> on <w:tblPr>).\n#    WdRowAlignment: wdAlignRowLeft=0, wdAlignRowCenter=1, wdAlignRowRight=2\n$t.Alignment = 1\n\n# 2) Round the \"Cot.M\u00e1xima R$\" column (3rd column) values.\n#    (Cell.Range.Text carries a trailing cell-mark, so the rounded values\n#    are written explicitly rather than parsed back out of the range.)\n$roundedValues = @(\"61.68\", \"61.67\", \"6.07\", \"19.5\", \"9.55\", \"6.13\")\n\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $cell = $t.Cell($i, 3)\n    $cell.Range.Text = $roundedValues[$i - 1]\n}\n"}
